$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '40.150.79'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '2.222.61'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.99'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.24%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("E13").Value = '  +3.84%  '
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '2.566.81'
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").Value = '2.234.65'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").Value = '40.074.83'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  +1.91%  '
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E28").Value = '  +3.47%  '
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("E30").Value = '  -5.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.06'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0716'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("E38").Value = '  +1.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.21%  '
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").Value = '2.071.75'
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.94%  '
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.79'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.87%  '
$ws.Range("E48").Value = '  -11.08%  '
$ws.Range("D49").Value = '2.438.54'
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("E50").Value = '  +4.80%  '
$ws.Range("E51").Value = '  +1.87%  '
